$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.269.77"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").Value = "2.799.10"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "116.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.552"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.10%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.18"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.13"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.85"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").Value = "3.240.63"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "2.779.92"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.895"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "52.187.85"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("E19").Value = "  +7.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.50%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.21"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.22"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("E25").Value = "  +5.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.73"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.25"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0827"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0416"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +17.51%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.93"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +22.11%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "128.19"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.53"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.115"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.31"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.074.36"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.983"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +15.98%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.98"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.91%  "
